$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift rows 121..155 down to 122..156 to make room for a new row 121 ---
# Work bottom-up so we don't overwrite data we still need to read.
for ($r = 155; $r -ge 121; $r--) {
    $ws.Range("A" + $r + ":C" + $r).Copy($ws.Range("A" + ($r + 1) + ":C" + ($r + 1)))
}

# --- Populate the new row 121 ---
# Borrow the formatting (font/fill/border/alignment) of row 3, which already
# carries the "A: style4, B: style6(shared text), C: style1" pattern used by
# the odd-numbered data rows.
$ws.Range("A3:C3").Copy()
$ws.Range("A121:C121").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A121").Value = 3017
$ws.Range("B121").Value = "플루리움4단지"
$ws.Range("C121").Value = 84

# New apartment-name font for this row (Dotum, matches the new 5th font).
$ws.Range("B121").Font.Name = "돋움"
$ws.Range("B121").Font.Size = 9

# Match the slightly shorter row height used for the new row.
$ws.Rows("121").RowHeight = 22.8

# --- Restore the view: scrolled near the new row, with I120:J121 selected ---
$excel.ActiveWindow.ScrollRow = 109
$ws.Range("I120:J121").Select()
